$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.464.97'
$ws.Range("E2").Value = '  +3.95%  '
$ws.Range("D3").Value = '3.368.55'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +3.88%  '
$ws.Range("E9").Value = '  +4.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.587'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000281'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.25%  '
$ws.Range("B13").Value = 'BitcoinCash'
$ws.Range("C13").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '642.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +12.54%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.909.67'
$ws.Range("E14").Value = '  +1.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.53%  '
$ws.Range("D16").Value = '68.569.50'
$ws.Range("E16").Value = '  +4.17%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.119'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.90%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.376.96'
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("E20").Value = '  +2.45%  '
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("E25").Value = '  +4.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '613.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.56%  '
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("D33").Value = '3.989.28'
$ws.Range("E33").Value = '  +6.74%  '
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("E35").Value = '  +2.87%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.41%  '
$ws.Range("E38").Value = '  +8.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.33'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '33.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").Value = '0.0₃0706'
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.344'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.21%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("E45").Value = '  +4.05%  '
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("E47").Value = '  +3.62%  '
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.29%  '
